$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet: Clientes -> Proveedores
$ws.Name = "Proveedores"

# Update the "Identificacion" example value and its example number (F4)
$ws.Range("E4").Value = "F"
$ws.Range("F4").Value = 111111111

# Replace the validation list for the identification-type column with the new letter codes
$ws.Range("E4").Validation.Formula1 = """F, J, D, N, E, O"""

# Remove the "Exento" and "EmisorReceptor" columns from the table definition first
# (so the table shrinks cleanly), then delete the now-unused worksheet columns.
$lo = $ws.ListObjects.Item(1)
$lo.ListColumns.Item(18).Delete()
$lo.ListColumns.Item(17).Delete()
$ws.Columns("Q:R").Delete()

# Update the comment on E2 describing the identification codes
$c = $ws.Range("E2").Comment
$newText = "Código de identificación:" + [char]10 + "F - Cédula física" + [char]10 + "J - Cédula jurídica" + [char]10 + "D - DIMEX" + [char]10 + "N - NITE" + [char]10 + "E - Extranjero" + [char]10 + "O - Otro"
$c.Text($newText)

# Move the active selection to A2
$ws.Range("A2").Select()

# Set the page to portrait orientation
$ws.PageSetup.Orientation = 1
